$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "r345"
$ws.Range("B8").Value = "rob"
$ws.Range("C8").Value = "is this in eastern now?"
$ws.Range("D8").Value = "2025-09-30 16:39:41"
